# Applies the weekly fruit/vegetable price-data refresh: the rows for
# "Hortaliza, Vega Modelo de Temuco - Locoto" (rows 2-15) keep their
# fixed descriptive columns (A,B,C,E,F,G,H,I,N,O,Q,R) but the per-record
# facts (Fecha, Volumen, Precio minimo/maximo/promedio, Precio $/Kg -
# columns D, J, K, L, M, P) get reshuffled across rows as the weekly
# source data is re-synced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in the refresh, by 1-based column index.
# D=4 (Fecha), J=10 (Volumen), K=11 (Precio minimo), L=12 (Precio maximo),
# M=13 (Precio promedio ponderado), P=16 (Precio $/Kg)
$cols = @(4, 10, 11, 12, 13, 16)

# Snapshot the current (pre-edit) values for every affected row/column
# before any writes happen, since several rows trade values with each
# other (cycles), not just simple overwrites.
$snapshot = @{}
foreach ($r in 2..15) {
    foreach ($c in $cols) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# Mapping: new row -> old row whose snapshot values it should receive.
$rowMap = @{
    2  = 9
    3  = 12
    4  = 4
    5  = 11
    6  = 14
    7  = 7
    8  = 8
    9  = 2
    10 = 10
    11 = 3
    12 = 15
    13 = 5
    14 = 6
    15 = 13
}

foreach ($newRow in 2..15) {
    $oldRow = $rowMap[$newRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($newRow, $c).Value = $snapshot["$oldRow-$c"]
    }
}
